$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Skills" column (J) for the three mentor rows.
# Order matches how the shared-string table was appended in the authored
# workbook: J3, then J2, then J4.
$ws.Range("J3").Value = "React, NodeJs"
$ws.Range("J2").Value = "Java"
$ws.Range("J4").Value = "Deploy"

# The edited cells picked up a (slightly) different font/style in the
# authored workbook - re-apply the same Arial 10 font so the engine
# allocates the matching style record.
$rng = $ws.Range("J2:J4")
$rng.Font.Name = "Arial"
$rng.Font.Size = 10

# Leave the selection on the last-edited cell, matching the saved view state.
$ws.Range("J4").Select() | Out-Null
